# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# to reflect newly generated counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of sheet name -> list of (row, newValue) updates for column F
$updates = @{
    "展览"   = @(
        @{ Row = 2;  Value = 6402 },
        @{ Row = 5;  Value = 377 },
        @{ Row = 7;  Value = 4 },
        @{ Row = 9;  Value = 81 },
        @{ Row = 10; Value = 72 },
        @{ Row = 14; Value = 776 },
        @{ Row = 15; Value = 3121 },
        @{ Row = 17; Value = 184 },
        @{ Row = 18; Value = 1786 }
    )
    "全部类型" = @(
        @{ Row = 2;  Value = 6402 },
        @{ Row = 5;  Value = 377 },
        @{ Row = 7;  Value = 4 },
        @{ Row = 10; Value = 81 },
        @{ Row = 11; Value = 72 },
        @{ Row = 15; Value = 776 },
        @{ Row = 16; Value = 3121 },
        @{ Row = 18; Value = 184 },
        @{ Row = 19; Value = 1786 }
    )
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($u in $updates[$sheetName]) {
        $ws.Cells.Item($u.Row, 6).Value = $u.Value
    }
}
